$d = $word.ActiveDocument

# Pull the whole package (all parts: document.xml, comments.xml, ...) as
# flattened WordprocessingML 2003 ("Flat OPC") text so we can strip the
# <w:contextualSpacing w:val="0"/> elements that were left over from the
# previous paragraph-properties pass wherever they occur.
$xml = $d.WordOpenXML

$xml = $xml -replace '<w:contextualSpacing[^>]*/>', ''

$d.WordOpenXML = $xml
